$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 331.66666
$ws.Range("I42").Value = 56.5
$ws.Range("K42").Value = 169.5
$ws.Range("M42").Value = 60.5
$ws.Range("H58").Value = 747.6875
$ws.Range("J58").Value = 1499.5
$ws.Range("L58").Value = 4498.5
$ws.Range("N58").Value = -4798.5
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 3000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3344
$ws.Range("H62").Value = 6143.1177
$ws.Range("I62").Value = 6162.1333
$ws.Range("K62").Value = 6162.1333
$ws.Range("M62").Value = -5538.1333
$ws.Range("H65").Value = 6143.1177
$ws.Range("I65").Value = 6162.1333
$ws.Range("K65").Value = 30810.6665
$ws.Range("M65").Value = -27690.6665
$ws.Range("H99").Value = 705.1111
$ws.Range("I99").Value = 431.16666
$ws.Range("J99").Value = 1253
$ws.Range("K99").Value = 1293.49998
$ws.Range("L99").Value = 3759
$ws.Range("M99").Value = 204.5000199999999
$ws.Range("N99").Value = -6755
$ws.Range("H104").Value = 571.1667
$ws.Range("I104").Value = 485.4
$ws.Range("K104").Value = 1456.2
$ws.Range("M104").Value = 290.8000000000002
$ws.Range("H115").Value = 31250580
$ws.Range("I115").Value = 740
$ws.Range("K115").Value = 2220
$ws.Range("M115").Value = -653
$ws.Range("H118").Value = 988.0769
$ws.Range("I118").Value = 804.1818
$ws.Range("K118").Value = 2412.5454
$ws.Range("M118").Value = -755.5454
$ws.Range("H127").Value = 1133.0714
$ws.Range("I127").Value = 1092.8334
$ws.Range("K127").Value = 3278.5002
$ws.Range("M127").Value = 1681.4998
$ws.Range("H129").Value = 2634.1538
$ws.Range("I129").Value = 1926.6154
$ws.Range("K129").Value = 5779.8462
$ws.Range("M129").Value = -779.8462
$ws.Range("H138").Value = 362763.47
$ws.Range("I138").Value = 7749.5
$ws.Range("J138").Value = 373205.06
$ws.Range("K138").Value = 23248.5
$ws.Range("L138").Value = 1119615.18
$ws.Range("M138").Value = -18108.5
$ws.Range("N138").Value = -1129895.18
$ws.Range("H139").Value = 69998.91
$ws.Range("J139").Value = 69998.91
$ws.Range("L139").Value = 69998.91
$ws.Range("N139").Value = -80278.91

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.4839
$ws.Range("I2").Value = 1514.1538
$ws.Range("J2").Value = 2037
$ws.Range("K2").Value = 1514.1538
$ws.Range("L2").Value = 2037
$ws.Range("M2").Value = -1401.1538
$ws.Range("N2").Value = -2263
$ws.Range("H31").Value = 2649.3333
$ws.Range("I31").Value = 1179.2
$ws.Range("K31").Value = 1179.2
$ws.Range("M31").Value = -885.2
$ws.Range("H32").Value = 14674.173
$ws.Range("I32").Value = 11133.703
$ws.Range("K32").Value = 11133.703
$ws.Range("M32").Value = -10846.703
$ws.Range("H102").Value = 5584.8276
$ws.Range("I102").Value = 5446.04
$ws.Range("K102").Value = 5446.04
$ws.Range("M102").Value = -3824.04
$ws.Range("H116").Value = 1598.4839
$ws.Range("I116").Value = 1514.1538
$ws.Range("J116").Value = 2037
$ws.Range("K116").Value = 1514.1538
$ws.Range("L116").Value = 2037
$ws.Range("M116").Value = 779.8462
$ws.Range("N116").Value = -6625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.4839
$ws.Range("I3").Value = 1514.1538
$ws.Range("J3").Value = 2037
$ws.Range("K3").Value = 1514.1538
$ws.Range("L3").Value = 2037
$ws.Range("M3").Value = -1400.1538
$ws.Range("N3").Value = -2265
$ws.Range("H105").Value = 13002679
$ws.Range("I105").Value = 716391.7
$ws.Range("K105").Value = 716391.7
$ws.Range("M105").Value = -714644.7
$ws.Range("H134").Value = 3790.9524
$ws.Range("I134").Value = 3307.3333
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 9921.999899999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -7386.999899999999
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8671.371999999999
$ws.Range("I31").Value = 6226.362
$ws.Range("J31").Value = 12378.322
$ws.Range("K31").Value = 6226.362
$ws.Range("L31").Value = 12378.322
$ws.Range("M31").Value = -5931.362
$ws.Range("N31").Value = -12968.322
$ws.Range("H34").Value = 8671.371999999999
$ws.Range("I34").Value = 6226.362
$ws.Range("J34").Value = 12378.322
$ws.Range("K34").Value = 6226.362
$ws.Range("L34").Value = 12378.322
$ws.Range("M34").Value = -6024.362
$ws.Range("N34").Value = -12782.322
$ws.Range("H99").Value = 14163.25
$ws.Range("I99").Value = 16835.334
$ws.Range("K99").Value = 16835.334
$ws.Range("M99").Value = -15337.334
$ws.Range("H126").Value = 14163.25
$ws.Range("I126").Value = 16835.334
$ws.Range("K126").Value = 50506.00199999999
$ws.Range("M126").Value = -48036.00199999999
$ws.Range("H134").Value = 3635.3547
$ws.Range("I134").Value = 2603.0908
$ws.Range("J134").Value = 6158.6665
$ws.Range("K134").Value = 7809.2724
$ws.Range("L134").Value = 18475.9995
$ws.Range("M134").Value = -5274.2724
$ws.Range("N134").Value = -23545.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 133.33333
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H103").Value = 20039.8
$ws.Range("I103").Value = 3377
$ws.Range("K103").Value = 10131
$ws.Range("M103").Value = -9252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2083.476
$ws.Range("I102").Value = 1718.2
$ws.Range("J102").Value = 2996.6667
$ws.Range("K102").Value = 1718.2
$ws.Range("L102").Value = 2996.6667
$ws.Range("M102").Value = -96.20000000000005
$ws.Range("N102").Value = -6240.6667
$ws.Range("H126").Value = 9104.823
$ws.Range("I126").Value = 5601.5
$ws.Range("J126").Value = 12218.889
$ws.Range("K126").Value = 16804.5
$ws.Range("L126").Value = 36656.667
$ws.Range("M126").Value = -14334.5
$ws.Range("N126").Value = -41596.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7424.9565
$ws.Range("I7").Value = 5769
$ws.Range("K7").Value = 5769
$ws.Range("M7").Value = -5657
$ws.Range("H40").Value = 9985.419
$ws.Range("I40").Value = 9850
$ws.Range("J40").Value = 9999.929
$ws.Range("K40").Value = 9850
$ws.Range("L40").Value = 9999.929
$ws.Range("M40").Value = -9714
$ws.Range("N40").Value = -10271.929
$ws.Range("H100").Value = 1596371.4
$ws.Range("I100").Value = 1936665.4
$ws.Range("J100").Value = 8333
$ws.Range("K100").Value = 1936665.4
$ws.Range("L100").Value = 8333
$ws.Range("M100").Value = -1936124.4
$ws.Range("N100").Value = -9415
$ws.Range("H126").Value = 7424.9565
$ws.Range("I126").Value = 5769
$ws.Range("K126").Value = 17307
$ws.Range("M126").Value = -14837
$ws.Range("H136").Value = 8764.85
$ws.Range("I136").Value = 6536.533
$ws.Range("K136").Value = 19609.599
$ws.Range("M136").Value = -17059.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3390.8
$ws.Range("I126").Value = 3306.2
$ws.Range("J126").Value = 3560
$ws.Range("K126").Value = 9918.599999999999
$ws.Range("L126").Value = 10680
$ws.Range("M126").Value = -7448.599999999999
$ws.Range("N126").Value = -15620
$ws.Range("H132").Value = 3370423.8
$ws.Range("I132").Value = 4068389.8
$ws.Range("K132").Value = 12205169.4
$ws.Range("M132").Value = -12202639.4
$ws.Range("H136").Value = 2643.875
$ws.Range("J136").Value = 5250
$ws.Range("L136").Value = 15750
$ws.Range("N136").Value = -20850
